# Trade #22 closed at 2026-02-17 12:37:12 - unknown UNKNOWN +0.000%
#
# This script updates the "Summary" and "Strategy Status" sheets with the
# refreshed aggregate statistics, and appends the newly closed trade (#22)
# as a new row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet - refresh aggregate metrics
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.59   # Current Capital
$summary.Range("B4").Value = 0.59      # Total P&L $
$summary.Range("B6").Value = 22        # Total Trades
$summary.Range("B7").Value = 8         # Winning Trades
$summary.Range("B9").Value = 36.36     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - refresh MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.59     # Capital
$status.Range("D4").Value = 22         # Trades
$status.Range("E4").Value = 0.59       # P&L $
$status.Range("F4").Value = 0.59       # P&L %
$status.Range("G4").Value = 36.36      # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade (row 23) to a worksheet in A:Q layout
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $row = 23

    $ws.Cells.Item($row, 1).Value = 22    # A: Trade #

    # Date / Time need to be written as plain text, not auto-converted
    # into Excel date/time serials.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"   # B: Date
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "12:37:06"     # C: Time
    $ws.Cells.Item($row, 3).ClearFormats()

    $ws.Cells.Item($row, 4).Value = "MarketMaking" # D: Strategy
    $ws.Cells.Item($row, 5).Value = "DOWN"         # E: Side
    $ws.Cells.Item($row, 6).Value = 0.41           # F: Entry Price
    $ws.Cells.Item($row, 7).Value = 0.43           # G: Exit Price
    $ws.Cells.Item($row, 8).Value = "CLOSED"       # H: Status
    $ws.Cells.Item($row, 9).Value = 4.878          # I: P&L %
    $ws.Cells.Item($row, 10).Value = 0.02          # J: P&L $
    $ws.Cells.Item($row, 11).Value = 100.59        # K: Capital After
    $ws.Cells.Item($row, 12).Value = 0             # L: Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0             # M: Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6           # N: Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps" # O: Entry Reason
    $ws.Cells.Item($row, 16).Value = "early_exit"  # P: Exit Reason
    $ws.Cells.Item($row, 17).Value = 0.15          # Q: Duration (min)
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
